$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing 2-column table; it will be recreated at the new A1:C15 extent.
$lo = $ws.ListObjects.Item(1)
$lo.Delete()

$data = @(
    @("Year", "Race/Ethnicity", "Percent Increase in Adjusted Rate"),
    @("2019 to 2020", "American Indian or Alaska Native", 0.202),
    @("2019 to 2020", "Asian", 0.219),
    @("2019 to 2020", "Black", 0.195),
    @("2019 to 2020", "Latino", 0.343),
    @("2019 to 2020", "Native Hawaiian and other Pacific Islander", 0.193),
    @("2019 to 2020", "Total", 0.158),
    @("2019 to 2020", "White", 0.076),
    @("2019 to 2021", "American Indian or Alaska Native", 0.38),
    @("2019 to 2021", "Asian", 0.22),
    @("2019 to 2021", "Black", 0.182),
    @("2019 to 2021", "Latino", 0.383),
    @("2019 to 2021", "Native Hawaiian and other Pacific Islander", 0.316),
    @("2019 to 2021", "Total", 0.175),
    @("2019 to 2021", "White", 0.082)
)

$r = 1
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Column widths: A (Year) narrower, B (Race/Ethnicity) same as old A, C (Percent...) new.
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 44
$ws.Columns.Item(3).ColumnWidth = 35

# Recreate the table over the new A1:C15 range with the same name/style.
$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:C15"), $null, 1)
$newLo.Name = "Table3"
$newLo.TableStyle = "TableStyleLight8"

Write-Host "done"
